$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1439.125
$ws.Range("I58").Value = 380.75
$ws.Range("J58").Value = 2497.5
$ws.Range("K58").Value = 1142.25
$ws.Range("L58").Value = 7492.5
$ws.Range("M58").Value = -992.25
$ws.Range("N58").Value = -7792.5
$ws.Range("H99").Value = 419.2
$ws.Range("I99").Value = 449.66666
$ws.Range("K99").Value = 1348.99998
$ws.Range("M99").Value = 149.0000199999999
$ws.Range("H100").Value = 3531.7273
$ws.Range("I100").Value = 3328.8572
$ws.Range("J100").Value = 3886.75
$ws.Range("K100").Value = 3328.8572
$ws.Range("L100").Value = 3886.75
$ws.Range("M100").Value = -2787.8572
$ws.Range("N100").Value = -4968.75
$ws.Range("H111").Value = 2938.5454
$ws.Range("J111").Value = 1257
$ws.Range("L111").Value = 3771
$ws.Range("N111").Value = -9905
$ws.Range("H118").Value = 516.7778
$ws.Range("I118").Value = 531.375
$ws.Range("K118").Value = 1594.125
$ws.Range("M118").Value = 62.875
$ws.Range("H127").Value = 3549
$ws.Range("I127").Value = 1300
$ws.Range("J127").Value = 4111.25
$ws.Range("K127").Value = 3900
$ws.Range("L127").Value = 12333.75
$ws.Range("M127").Value = 1060
$ws.Range("N127").Value = -22253.75
$ws.Range("H129").Value = 13946.471
$ws.Range("J129").Value = 25079.555
$ws.Range("L129").Value = 75238.66500000001
$ws.Range("N129").Value = -85238.66500000001
$ws.Range("H137").Value = 12791.615
$ws.Range("I137").Value = 6558.143
$ws.Range("K137").Value = 19674.429
$ws.Range("M137").Value = -17124.429

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4803.8936
$ws.Range("I32").Value = 4292.524
$ws.Range("J32").Value = 9099.4
$ws.Range("K32").Value = 4292.524
$ws.Range("L32").Value = 9099.4
$ws.Range("M32").Value = -4005.524
$ws.Range("N32").Value = -9673.4
$ws.Range("H61").Value = 27779958
$ws.Range("I61").Value = 31251952
$ws.Range("K61").Value = 31251952
$ws.Range("M61").Value = -31251740
$ws.Range("H74").Value = 45459308
$ws.Range("I74").Value = 52636504
$ws.Range("K74").Value = 52636504
$ws.Range("M74").Value = -52635630
$ws.Range("H77").Value = 45459308
$ws.Range("I77").Value = 52636504
$ws.Range("K77").Value = 263182520
$ws.Range("M77").Value = -263178152
$ws.Range("H132").Value = 34542744
$ws.Range("I132").Value = 14551.137
$ws.Range("J132").Value = 143059940
$ws.Range("K132").Value = 43653.411
$ws.Range("L132").Value = 429179820
$ws.Range("M132").Value = -41123.411
$ws.Range("N132").Value = -429184880
$ws.Range("H136").Value = 27779958
$ws.Range("I136").Value = 31251952
$ws.Range("K136").Value = 93755856
$ws.Range("M136").Value = -93753306

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1985.0714
$ws.Range("I94").Value = 1225.75
$ws.Range("J94").Value = 2997.5
$ws.Range("K94").Value = 1225.75
$ws.Range("L94").Value = 2997.5
$ws.Range("M94").Value = -774.75
$ws.Range("N94").Value = -3899.5
$ws.Range("H134").Value = 5625
$ws.Range("I134").Value = 5625
$ws.Range("K134").Value = 16875
$ws.Range("M134").Value = -14340

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1705.4445
$ws.Range("I16").Value = 1731.125
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1731.125
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1444.125
$ws.Range("N16").Value = -2074
$ws.Range("H22").Value = 1607.9166
$ws.Range("I22").Value = 436.875
$ws.Range("K22").Value = 436.875
$ws.Range("M22").Value = -86.875
$ws.Range("H113").Value = 1705.4445
$ws.Range("I113").Value = 1731.125
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1731.125
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 438.875
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 67519.25999999999
$ws.Range("I132").Value = 71845.14
$ws.Range("K132").Value = 215535.42
$ws.Range("M132").Value = -213005.42
$ws.Range("H134").Value = 2843.8276
$ws.Range("I134").Value = 2350.348
$ws.Range("K134").Value = 7051.044
$ws.Range("M134").Value = -4516.044

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 235
$ws.Range("I10").Value = 300
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 120
$ws.Range("M10").Value = -761
$ws.Range("N10").Value = -398
$ws.Range("H26").Value = 460.75
$ws.Range("I26").Value = 258.5
$ws.Range("K26").Value = 775.5
$ws.Range("M26").Value = -487.5
$ws.Range("H36").Value = 1619
$ws.Range("I36").Value = 586
$ws.Range("J36").Value = 2032.2
$ws.Range("K36").Value = 1758
$ws.Range("L36").Value = 6096.6
$ws.Range("M36").Value = -1589
$ws.Range("N36").Value = -6434.6
$ws.Range("H115").Value = 7057.5
$ws.Range("J115").Value = 9326.666999999999
$ws.Range("L115").Value = 27980.001
$ws.Range("N115").Value = -30330.001
$ws.Range("H121").Value = 1217.7142
$ws.Range("J121").Value = 841.6667
$ws.Range("L121").Value = 2525.0001
$ws.Range("N121").Value = -5145.0001
$ws.Range("H137").Value = 3910.4
$ws.Range("I137").Value = 2607.6
$ws.Range("J137").Value = 5213.2
$ws.Range("K137").Value = 7822.799999999999
$ws.Range("L137").Value = 15639.6
$ws.Range("M137").Value = -2722.799999999999
$ws.Range("N137").Value = -25839.6
$ws.Range("H140").Value = 2014.6923
$ws.Range("I140").Value = 1612.5
$ws.Range("K140").Value = 4837.5
$ws.Range("M140").Value = 342.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 760004.75
$ws.Range("J33").Value = 1006673
$ws.Range("L33").Value = 1006673
$ws.Range("N33").Value = -1007177
$ws.Range("H41").Value = 7933.3335
$ws.Range("I41").Value = 7933.3335
$ws.Range("K41").Value = 7933.3335
$ws.Range("M41").Value = -7578.3335
$ws.Range("H63").Value = 39995
$ws.Range("J63").Value = 39995
$ws.Range("L63").Value = 39995
$ws.Range("N63").Value = -41367
$ws.Range("H66").Value = 39995
$ws.Range("J66").Value = 39995
$ws.Range("L66").Value = 119985
$ws.Range("N66").Value = -126849

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3799.8
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3799.8
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3799.8
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -4389.8
$ws.Range("H27").Value = 3799.8
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3799.8
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3799.8
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4013.8
$ws.Range("H61").Value = 3465.2
$ws.Range("I61").Value = 2921.4614
$ws.Range("J61").Value = 6999.5
$ws.Range("K61").Value = 2921.4614
$ws.Range("L61").Value = 6999.5
$ws.Range("M61").Value = -2719.4614
$ws.Range("N61").Value = -7403.5
$ws.Range("H113").Value = 3465.2
$ws.Range("I113").Value = 2921.4614
$ws.Range("J113").Value = 6999.5
$ws.Range("K113").Value = 2921.4614
$ws.Range("L113").Value = 6999.5
$ws.Range("M113").Value = -751.4614000000001
$ws.Range("N113").Value = -11339.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 47497.5
$ws.Range("H67").Value = 47497.5
$ws.Range("H107").Value = 1705.1333
$ws.Range("I107").Value = 1234.2727
$ws.Range("K107").Value = 3702.8181
$ws.Range("M107").Value = -1782.8181
$ws.Range("H113").Value = 860.5
$ws.Range("I113").Value = 861.8823
$ws.Range("J113").Value = 857.8889
$ws.Range("K113").Value = 2585.6469
$ws.Range("L113").Value = 2573.6667
$ws.Range("M113").Value = -415.6468999999997
$ws.Range("N113").Value = -6913.6667
